$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Datos actualizados a 8 de Abril de 2020 a las 12:22'

$ws.Range('A10').Value = 'Iran'
$ws.Range('B10').Value = 67286
$ws.Range('C10').Value = 4697
$ws.Range('D10').Value = 27039
$ws.Range('E10').Value = 36244
$ws.Range('F10').Value = 3987
$ws.Range('G10').Value = 131
$ws.Range('H10').Value = 4003

$ws.Range('A14').Value = 'Suiza'
$ws.Range('B14').Value = 22488
$ws.Range('C14').Value = 235
$ws.Range('D14').Value = 8704
$ws.Range('E14').Value = 12938
$ws.Range('F14').Value = 391
$ws.Range('G14').Value = 25
$ws.Range('H14').Value = 846

$ws.Range('A18').Value = 'Austria'
$ws.Range('B18').Value = 12738
$ws.Range('C18').Value = 99
$ws.Range('D18').Value = 4512
$ws.Range('E18').Value = 7953
$ws.Range('F18').Value = 267
$ws.Range('G18').Value = 30
$ws.Range('H18').Value = 273

$ws.Range('A32').Value = 'Rumania'
$ws.Range('B32').Value = 4761
$ws.Range('C32').Value = 344
$ws.Range('D32').Value = 528
$ws.Range('E32').Value = 4024
$ws.Range('F32').Value = 162
$ws.Range('G32').Value = 12
$ws.Range('H32').Value = 209

$ws.Range('A60').Value = 'Marruecos'
$ws.Range('B60').Value = 1242
$ws.Range('C60').Value = 58
$ws.Range('D60').Value = 97
$ws.Range('E60').Value = 1054
$ws.Range('F60').Value = 1
$ws.Range('G60').Value = 1
$ws.Range('H60').Value = 91

$ws.Range('A61').Value = 'Nueva Zelanda'
$ws.Range('B61').Value = 1210
$ws.Range('C61').Value = 50
$ws.Range('D61').Value = 282
$ws.Range('E61').Value = 927
$ws.Range('F61').Value = 4
$ws.Range('G61').Value = 0
$ws.Range('H61').Value = 1

$ws.Range('A62').Value = 'Estonia'
$ws.Range('B62').Value = 1185
$ws.Range('C62').Value = 36
$ws.Range('D62').Value = 72
$ws.Range('E62').Value = 1089
$ws.Range('F62').Value = 11
$ws.Range('G62').Value = 3
$ws.Range('H62').Value = 24

$ws.Range('A83').Value = 'Libano'
$ws.Range('B83').Value = 575
$ws.Range('C83').Value = 27
$ws.Range('D83').Value = 62
$ws.Range('E83').Value = 494
$ws.Range('F83').Value = 29
$ws.Range('G83').Value = 0
$ws.Range('H83').Value = 19

$ws.Range('A91').Value = 'Albania'
$ws.Range('B91').Value = 400
$ws.Range('C91').Value = 17
$ws.Range('D91').Value = 154
$ws.Range('E91').Value = 224
$ws.Range('F91').Value = 7
$ws.Range('G91').Value = 0
$ws.Range('H91').Value = 22

$ws.Range('A92').Value = 'Cuba'
$ws.Range('B92').Value = 396
$ws.Range('C92').Value = 0
$ws.Range('D92').Value = 27
$ws.Range('E92').Value = 358
$ws.Range('F92').Value = 15
$ws.Range('G92').Value = 0
$ws.Range('H92').Value = 11

$ws.Range('A93').Value = 'Burkina Faso'
$ws.Range('B93').Value = 384
$ws.Range('C93').Value = 0
$ws.Range('D93').Value = 127
$ws.Range('E93').Value = 238
$ws.Range('F93').Value = 0
$ws.Range('G93').Value = 0
$ws.Range('H93').Value = 19

$ws.Range('A110').Value = 'Banglades'
$ws.Range('B110').Value = 218
$ws.Range('C110').Value = 54
$ws.Range('D110').Value = 33
$ws.Range('E110').Value = 165
$ws.Range('F110').Value = 1
$ws.Range('G110').Value = 3
$ws.Range('H110').Value = 20

$ws.Range('A111').Value = 'Bolivia'
$ws.Range('B111').Value = 210
$ws.Range('C111').Value = 16
$ws.Range('D111').Value = 2
$ws.Range('E111').Value = 193
$ws.Range('F111').Value = 3
$ws.Range('G111').Value = 1
$ws.Range('H111').Value = 15

$ws.Range('A112').Value = 'Georgia'
$ws.Range('B112').Value = 208
$ws.Range('C112').Value = 12
$ws.Range('D112').Value = 48
$ws.Range('E112').Value = 157
$ws.Range('F112').Value = 6
$ws.Range('G112').Value = 0
$ws.Range('H112').Value = 3

$ws.Range('A113').Value = 'Sri Lanka'
$ws.Range('B113').Value = 186
$ws.Range('C113').Value = 1
$ws.Range('D113').Value = 42
$ws.Range('E113').Value = 138
$ws.Range('F113').Value = 5
$ws.Range('G113').Value = 0
$ws.Range('H113').Value = 6

$ws.Range('A114').Value = 'Islas Feroe'
$ws.Range('B114').Value = 184
$ws.Range('C114').Value = 0
$ws.Range('D114').Value = 131
$ws.Range('E114').Value = 53
$ws.Range('F114').Value = 1
$ws.Range('G114').Value = 0
$ws.Range('H114').Value = 0

$ws.Range('A115').Value = 'Consejo Danes para los Refugiados'
$ws.Range('B115').Value = 180
$ws.Range('C115').Value = 0
$ws.Range('D115').Value = 9
$ws.Range('E115').Value = 153
$ws.Range('F115').Value = 0
$ws.Range('G115').Value = 0
$ws.Range('H115').Value = 18

$ws.Range('A116').Value = 'Kenia'
$ws.Range('B116').Value = 172
$ws.Range('C116').Value = 0
$ws.Range('D116').Value = 7
$ws.Range('E116').Value = 159
$ws.Range('F116').Value = 2
$ws.Range('G116').Value = 0
$ws.Range('H116').Value = 6

$ws.Range('A117').Value = 'Mayotte'
$ws.Range('B117').Value = 171
$ws.Range('C117').Value = 0
$ws.Range('D117').Value = 22
$ws.Range('E117').Value = 147
$ws.Range('F117').Value = 3
$ws.Range('G117').Value = 0
$ws.Range('H117').Value = 2

$ws.Range('A118').Value = 'Venezuela'
$ws.Range('B118').Value = 166
$ws.Range('C118').Value = 0
$ws.Range('D118').Value = 65
$ws.Range('E118').Value = 94
$ws.Range('F118').Value = 6
$ws.Range('G118').Value = 0
$ws.Range('H118').Value = 7

$ws.Range('A124').Value = 'Republica de Yibuti'
$ws.Range('B124').Value = 121
$ws.Range('C124').Value = 31
$ws.Range('D124').Value = 18
$ws.Range('E124').Value = 103
$ws.Range('F124').Value = 0
$ws.Range('G124').Value = 0
$ws.Range('H124').Value = 0

$ws.Range('A125').Value = 'Paraguay'
$ws.Range('B125').Value = 119
$ws.Range('C125').Value = 4
$ws.Range('D125').Value = 15
$ws.Range('E125').Value = 99
$ws.Range('F125').Value = 1
$ws.Range('G125').Value = 0
$ws.Range('H125').Value = 5

$ws.Range('A126').Value = 'Camboya'
$ws.Range('B126').Value = 117
$ws.Range('C126').Value = 2
$ws.Range('D126').Value = 63
$ws.Range('E126').Value = 54
$ws.Range('F126').Value = 1
$ws.Range('G126').Value = 0
$ws.Range('H126').Value = 0

$ws.Range('A127').Value = 'Gibraltar'
$ws.Range('B127').Value = 113
$ws.Range('C127').Value = 0
$ws.Range('D127').Value = 60
$ws.Range('E127').Value = 53
$ws.Range('F127').Value = 0
$ws.Range('G127').Value = 0
$ws.Range('H127').Value = 0

$ws.Range('A128').Value = 'Trinidad yTobago'
$ws.Range('B128').Value = 107
$ws.Range('C128').Value = 0
$ws.Range('D128').Value = 1
$ws.Range('E128').Value = 98
$ws.Range('F128').Value = 0
$ws.Range('G128').Value = 0
$ws.Range('H128').Value = 8

$ws.Range('A129').Value = 'Ruanda'
$ws.Range('B129').Value = 105
$ws.Range('C129').Value = 0
$ws.Range('D129').Value = 7
$ws.Range('E129').Value = 98
$ws.Range('F129').Value = 0
$ws.Range('G129').Value = 0
$ws.Range('H129').Value = 0

$ws.Range('A130').Value = 'El Salvador'
$ws.Range('B130').Value = 93
$ws.Range('C130').Value = 15
$ws.Range('D130').Value = 9
$ws.Range('E130').Value = 79
$ws.Range('F130').Value = 2
$ws.Range('G130').Value = 1
$ws.Range('H130').Value = 5

